# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Fri Nov 29 19:45:13 UTC 2024 with GitHub Actions".
#
# Numeric-looking text values (e.g. "656.26") must stay stored as TEXT
# (matching the original inline-string cells), so for those we briefly
# force NumberFormat="@" while assigning, then restore the cell style to
# "Normal" so no stray numeric style is left behind in styles.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "97.322.38"
$ws.Range("E2").Value = "  +2.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.580.11"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - Solana
$ws.Range("E5").Value = "  +2.54%  "

# Row 6 - BNB (was XRP)
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "656.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "

# Row 7 - XRP (was BNB)
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +17.25%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.11%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.06%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +5.14%  "

# Row 11 - LidoStakedEther
$ws.Range("D11").Value = "3.577.45"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.46%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.49%  "

# Row 14 - Toncoin
$ws.Range("E14").Value = "  -0.71%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.245.29"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "97.126.80"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000261"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.94%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.573.31"
$ws.Range("E19").Value = "  -0.08%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "

# Row 22 - Stellar
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.532"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.31%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +0.99%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "514.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +5.45%  "

# Row 26 - NEARProtocol
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "

# Row 27 - Litecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.12%  "

# Row 28 - Aptos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "3.772.60"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30 - Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.98%  "

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.56%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.15%  "

# Row 33 - Dai
$ws.Range("E33").Value = "  -0.02%  "

# Row 34 - Cronos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.185"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.53%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.59%  "

# Row 36 - EthereumClassic
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.12%  "

# Row 38 - Bittensor
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "615.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.71%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("E39").Value = "  +1.49%  "

# Row 40 - Fetch.AI
$ws.Range("E40").Value = "  -2.38%  "

# Row 41 - Kaspa (was ImmutableX)
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.155"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.79%  "

# Row 42 - ImmutableX (was Kaspa)
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.98%  "

# Row 44 - ARBITRUM
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.925"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45 - Filecoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.29%  "

# Row 46 - VeChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0440"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.15%  "

# Row 47 - Stacks
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.87%  "

# Row 48 - WhiteBITCoin
$ws.Range("E48").Value = "  +0.95%  "

# Row 49 - Algorand
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.412"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +34.11%  "

# Row 50 - Cosmos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.47%  "

# Row 51 - dogwifhat (was EnergySwap)
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.95%  "
